$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    3  = -4
    4  = 1
    5  = 5
    6  = -3
    7  = 0
    8  = 1
    9  = -2
    10 = 2
    11 = -3
    12 = 1
    13 = -1
    14 = 2
    15 = 3
    16 = 2
    17 = -3
    18 = -4
    19 = 3
    20 = -1
    21 = 1
    22 = 5
    23 = 2
    24 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
